$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue "D2" "65.336.33"
$ws.Range("E2").Value = "  -4.01%  "
Set-TextValue "D3" "3.390.87"
$ws.Range("E3").Value = "  -5.91%  "
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.19%  "
Set-TextValue "D5" "182.21"
$ws.Range("E5").Value = "  -10.92%  "
Set-TextValue "D6" "527.60"
$ws.Range("E6").Value = "  -6.60%  "
Set-TextValue "D7" "0.609"
$ws.Range("E7").Value = "  -2.18%  "
Set-TextValue "D8" "3.393.21"
$ws.Range("E8").Value = "  -5.65%  "
Set-TextValue "D9" "0.999"
$ws.Range("E9").Value = "  -0.16%  "
Set-TextValue "D10" "0.626"
$ws.Range("E10").Value = "  -7.00%  "
Set-TextValue "D11" "57.56"
$ws.Range("E11").Value = "  -6.24%  "
Set-TextValue "D12" "0.134"
$ws.Range("E12").Value = "  -11.69%  "
Set-TextValue "D13" "0.0000256"
$ws.Range("E13").Value = "  -10.74%  "
Set-TextValue "D14" "9.32"
$ws.Range("E14").Value = "  -7.13%  "
Set-TextValue "D15" "3.928.56"
$ws.Range("E15").Value = "  -6.29%  "
$ws.Range("E16").Value = "  -2.71%  "
Set-TextValue "D17" "3.382.78"
$ws.Range("E17").Value = "  -6.25%  "
Set-TextValue "D18" "64.959.70"
$ws.Range("E18").Value = "  -4.34%  "
Set-TextValue "D19" "17.53"
$ws.Range("E19").Value = "  -7.41%  "
Set-TextValue "D20" "11.23"
$ws.Range("E20").Value = "  -9.23%  "
Set-TextValue "D21" "0.975"
$ws.Range("E21").Value = "  -9.34%  "
Set-TextValue "D22" "374.75"
$ws.Range("E22").Value = "  -6.86%  "
Set-TextValue "D23" "82.78"
$ws.Range("E23").Value = "  -3.01%  "
Set-TextValue "D24" "3.72"
$ws.Range("E24").Value = "  -10.64%  "
Set-TextValue "D25" "10.85"
$ws.Range("E25").Value = "  -17.87%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D26" "3.63"
$ws.Range("E26").Value = "  -8.28%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D27" "11.54"
$ws.Range("E27").Value = "  -8.49%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D28" "2.66"
$ws.Range("E28").Value = "  -9.30%  "
Set-TextValue "D29" "8.51"
$ws.Range("E29").Value = "  -9.51%  "
Set-TextValue "D30" "676.90"
$ws.Range("E30").Value = "  +0.55%  "
Set-TextValue "D31" "29.70"
$ws.Range("E31").Value = "  -6.08%  "
Set-TextValue "D32" "6.74"
$ws.Range("E32").Value = "  -18.92%  "
Set-TextValue "D33" "61.64"
$ws.Range("E33").Value = "  -3.27%  "
Set-TextValue "D34" "11.16"
$ws.Range("E34").Value = "  -8.70%  "
Set-TextValue "D35" "0.106"
$ws.Range("E35").Value = "  -7.33%  "
$ws.Range("E36").Value = "  +0.20%  "
Set-TextValue "D37" "36.67"
$ws.Range("E37").Value = "  -13.34%  "
Set-TextValue "D38" "0.386"
$ws.Range("E38").Value = "  -8.27%  "
Set-TextValue "D39" "0.997"
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("E40").Value = "  -6.44%  "
Set-TextValue "D41" "2.869.59"
$ws.Range("E41").Value = "  -13.26%  "
Set-TextValue "D42" "2.77"
$ws.Range("E42").Value = "  -12.69%  "
Set-TextValue "D43" "2.64"
$ws.Range("E43").Value = "  -4.00%  "
Set-TextValue "D44" "0.0₃0627"
$ws.Range("E44").Value = "  -18.50%  "
$ws.Range("E45").Value = "  -7.42%  "
Set-TextValue "D46" "2.33"
$ws.Range("E46").Value = "  -15.89%  "
Set-TextValue "D47" "0.126"
$ws.Range("E47").Value = "  -4.47%  "
Set-TextValue "D48" "135.29"
$ws.Range("E48").Value = "  -3.23%  "
Set-TextValue "D49" "2.85"
$ws.Range("E49").Value = "  -7.09%  "
Set-TextValue "D50" "2.54"
$ws.Range("E50").Value = "  -6.97%  "
Set-TextValue "D51" "7.61"
$ws.Range("E51").Value = "  -14.00%  "
